$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44377
$ws.Range("I2").Value = "Segunda"
$ws.Range("J2").Value = 550
$ws.Range("K2").Value = 2000
$ws.Range("L2").Value = 2800
$ws.Range("M2").Value = 2364
$ws.Range("P2").Value = 394

# Row 3
$ws.Range("D3").Value = 44370
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 1200
$ws.Range("M3").Value = 1080
$ws.Range("P3").Value = 180

# Row 4
$ws.Range("D4").Value = 44267
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 1500
$ws.Range("L4").Value = 1800
$ws.Range("M4").Value = 1650
$ws.Range("P4").Value = 275
